$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New items appended to the end of the single-column item list (column A)
$items = @(
    "HRT Roll(Tissue roll)",
    "Banana",
    "Electric Material",
    "Plumbing Material",
    "Hardware material",
    "Kitchenware",
    "CP Pomegranate",
    "CP Watermelone",
    "CP Pineapple",
    "CP Blackcane",
    "CP Classic Cane",
    "CP ABC",
    "Incense Stick (Agarbatthi)"
)

$startRow = 296
$row = $startRow
foreach ($item in $items) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $item
    # Match the formatting already used by the rows directly above (thin
    # left/right borders, no top/bottom) so the new rows look consistent
    # with the rest of the list.
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
    $row = $row + 1
}

$lastRow = $row - 1

# Put the selection/active cell on the last new row, as Excel would leave
# it after typing the final entry.
$ws.Cells.Item($lastRow, 1).Select()
